$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A (timestamp) and B (label) are untouched by this edit - only the
# sensor reading columns C:H change value. Rows 2 and 3 receive two brand
# new samples, the data that used to live in rows 2-19 is copied down two
# rows into 4-21, and the old trailing row (22) is removed so the sheet ends
# up with one fewer data row (new dimension A1:H21).

$ws.Range("C2").Value = -0.03897037506103547
$ws.Range("D2").Value = -0.04966262578964246
$ws.Range("E2").Value = 0.6067050054669385
$ws.Range("F2").Value = -0.0704022198915481
$ws.Range("G2").Value = 0.1944078654050827
$ws.Range("H2").Value = -0.0245873257517814

$ws.Range("C3").Value = -0.3774656057357791
$ws.Range("D3").Value = -0.08384630084037778
$ws.Range("E3").Value = 1.005563378334046
$ws.Range("F3").Value = 0.3005456924438476
$ws.Range("G3").Value = 0.8894197940826416
$ws.Range("H3").Value = 0.086895577609539

$ws.Range("C4").Value = -0.389985084533691
$ws.Range("D4").Value = -0.08737320899963376
$ws.Range("E4").Value = 1.122915458679199
$ws.Range("F4").Value = 0.266642689704895
$ws.Range("G4").Value = 1.489289522171021
$ws.Range("H4").Value = 1.197143197059631

$ws.Range("C5").Value = -3.129088830947885
$ws.Range("D5").Value = -0.5355502605438246
$ws.Range("E5").Value = 2.341251343488698
$ws.Range("F5").Value = -0.5832235813140869
$ws.Range("G5").Value = 0.5451972484588623
$ws.Range("H5").Value = 0.808174729347229

$ws.Range("C6").Value = -4.376520133018492
$ws.Range("D6").Value = -0.8961358070373535
$ws.Range("E6").Value = 3.528440594673157
$ws.Range("F6").Value = -0.5074763298034668
$ws.Range("G6").Value = -0.9810495972633362
$ws.Range("H6").Value = 0.5829181671142578

$ws.Range("C7").Value = -2.786781752109524
$ws.Range("D7").Value = -1.147766584157945
$ws.Range("E7").Value = 5.211621630191809
$ws.Range("F7").Value = -0.9239336848258972
$ws.Range("G7").Value = -0.2371684312820434
$ws.Range("H7").Value = 0.1962404549121856

$ws.Range("C8").Value = 1.074200153350837
$ws.Range("D8").Value = 1.294980049133306
$ws.Range("E8").Value = 4.674310564994808
$ws.Range("F8").Value = -0.8011497855186462
$ws.Range("G8").Value = -1.048092007637024
$ws.Range("H8").Value = -0.0548251569271087

$ws.Range("C9").Value = 33.04098894596114
$ws.Range("D9").Value = -6.460391509532979
$ws.Range("E9").Value = -13.8664929449559
$ws.Range("F9").Value = 0.3634648323059082
$ws.Range("G9").Value = 1.29381263256073
$ws.Range("H9").Value = 1.935221076011657

$ws.Range("C10").Value = 33.17409253120405
$ws.Range("D10").Value = -12.13667659759519
$ws.Range("E10").Value = -16.41260833740226
$ws.Range("F10").Value = -4.552626132965088
$ws.Range("G10").Value = -3.393509149551392
$ws.Range("H10").Value = 1.527773976325989

$ws.Range("C11").Value = -4.199231290817252
$ws.Range("D11").Value = -1.979123908281309
$ws.Range("E11").Value = 2.645077538490293
$ws.Range("F11").Value = 0.7554876208305359
$ws.Range("G11").Value = 1.477530360221863
$ws.Range("H11").Value = -0.5998696684837341

$ws.Range("C12").Value = -1.133084297180169
$ws.Range("D12").Value = 1.222284126281733
$ws.Range("E12").Value = 2.30235185623169
$ws.Range("F12").Value = 0.3274237811565399
$ws.Range("G12").Value = -0.0329867228865623
$ws.Range("H12").Value = -0.4276056587696075

$ws.Range("C13").Value = 1.011236310005188
$ws.Range("D13").Value = -0.08591727167367935
$ws.Range("E13").Value = 1.941863000392914
$ws.Range("F13").Value = 0.0789543315768241
$ws.Range("G13").Value = -0.0847575515508651
$ws.Range("H13").Value = -0.1519527286291122

$ws.Range("C14").Value = 1.760050582885739
$ws.Range("D14").Value = -0.07424210608005596
$ws.Range("E14").Value = 0.6869683876633619
$ws.Range("F14").Value = -0.0221438650041818
$ws.Range("G14").Value = 0.1513418704271316
$ws.Range("H14").Value = -0.171500414609909

$ws.Range("C15").Value = -0.4296665787696987
$ws.Range("D15").Value = -0.05021018907427618
$ws.Range("E15").Value = 0.4772178567945989
$ws.Range("F15").Value = 0.0762054398655891
$ws.Range("G15").Value = -0.0978911519050598
$ws.Range("H15").Value = 0.366213709115982

$ws.Range("C16").Value = -1.509320116043069
$ws.Range("D16").Value = -0.4450684934854586
$ws.Range("E16").Value = 1.374938857555391
$ws.Range("F16").Value = -0.2686280012130737
$ws.Range("G16").Value = 1.315498352050781
$ws.Range("H16").Value = -0.0265726372599601

$ws.Range("C17").Value = 2.49392051696776
$ws.Range("D17").Value = -1.748301430046552
$ws.Range("E17").Value = 0.9690718531608501
$ws.Range("F17").Value = -0.1499674171209335
$ws.Range("G17").Value = -0.5490151643753052
$ws.Range("H17").Value = 0.4234823286533355

$ws.Range("C18").Value = -0.5127081871032715
$ws.Range("D18").Value = -0.141617327928543
$ws.Range("E18").Value = -0.73384278640151
$ws.Range("F18").Value = -0.09239336848258969
$ws.Range("G18").Value = -0.09071348607540131
$ws.Range("H18").Value = 0.0775798857212066

$ws.Range("C19").Value = -0.1017783880233756
$ws.Range("D19").Value = 0.3632039599120608
$ws.Range("E19").Value = 0.1947979252785462
$ws.Range("F19").Value = 0.0155770638957619
$ws.Range("G19").Value = -0.1032362282276153
$ws.Range("H19").Value = 0.09071348607540131

$ws.Range("C20").Value = 0.2760831832885773
$ws.Range("D20").Value = 0.05902776718139471
$ws.Range("E20").Value = 0.4707315444946266
$ws.Range("F20").Value = 0.4306600093841553
$ws.Range("G20").Value = 0.7470881938934326
$ws.Range("H20").Value = -0.1643227487802505

$ws.Range("C21").Value = 0.7326052427291816
$ws.Range("D21").Value = -0.04417074620723105
$ws.Range("E21").Value = 0.2964785575866773
$ws.Range("F21").Value = 1.219439744949341
$ws.Range("G21").Value = 1.55587375164032
$ws.Range("H21").Value = 0.009468411095440299

# Drop the now-superfluous last row so the sheet matches dimension A1:H21.
$ws.Rows.Item(22).Delete()
